# Scheduled-runner update: refresh cached market-board price snapshots
# (currentAveragePrice / NQ / HQ / Leve prices / profits) across the
# per-job Leve tables. Values only - no structural/formula changes.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 11 (ALC)
$ws.Range("H11").Value = 2700
$ws.Range("I11").Value = 2700
$ws.Range("K11").Value = 2700
$ws.Range("M11").Value = -2560

# Row 64 (ALC)
$ws.Range("H64").Value = 6981.4136
$ws.Range("I64").Value = 3348.25
$ws.Range("J64").Value = 11453
$ws.Range("K64").Value = 3348.25
$ws.Range("L64").Value = 11453
$ws.Range("M64").Value = -3100.25
$ws.Range("N64").Value = -11949

# Row 67 (ALC)
$ws.Range("H67").Value = 6981.4136
$ws.Range("I67").Value = 3348.25
$ws.Range("J67").Value = 11453
$ws.Range("K67").Value = 3348.25
$ws.Range("L67").Value = 11453
$ws.Range("M67").Value = -2490.25
$ws.Range("N67").Value = -13169

# Row 125 (ALC)
$ws.Range("H125").Value = 26565.2
$ws.Range("J125").Value = 997
$ws.Range("L125").Value = 8973
$ws.Range("N125").Value = -13893

# Row 132 (ALC)
$ws.Range("H132").Value = 1687.8
$ws.Range("I132").Value = 1479.7742
$ws.Range("K132").Value = 4439.3226
$ws.Range("M132").Value = -1909.3226

# Row 133 (ALC)
$ws.Range("H133").Value = 69749.5
$ws.Range("J133").Value = 69749.5
$ws.Range("L133").Value = 69749.5
$ws.Range("N133").Value = -79869.5

$ws = $wb.Worksheets.Item("ARM")
# Row 45 (ARM)
$ws.Range("H45").Value = 4840.8184
$ws.Range("I45").Value = 2750
$ws.Range("K45").Value = 2750
$ws.Range("M45").Value = -2373

# Row 97 (ARM)
$ws.Range("H97").Value = 1148.2307
$ws.Range("I97").Value = 836.5833
$ws.Range("K97").Value = 836.5833
$ws.Range("M97").Value = -340.5833

# Row 110 (ARM)
$ws.Range("H110").Value = 4423.0938
$ws.Range("I110").Value = 3232.625
$ws.Range("J110").Value = 7994.5
$ws.Range("K110").Value = 3232.625
$ws.Range("L110").Value = 7994.5
$ws.Range("M110").Value = -1187.625
$ws.Range("N110").Value = -12084.5

# Row 132 (ARM)
$ws.Range("H132").Value = 4899.8
$ws.Range("I132").Value = 5454.273
$ws.Range("K132").Value = 16362.819
$ws.Range("M132").Value = -13832.819

$ws = $wb.Worksheets.Item("BSM")
# Row 2 (BSM)
$ws.Range("H2").Value = 69149.5
$ws.Range("J2").Value = 69149.5
$ws.Range("L2").Value = 69149.5
$ws.Range("N2").Value = -69375.5

# Row 20 (BSM)
$ws.Range("H20").Value = 4324.0625
$ws.Range("I20").Value = 3811.158
$ws.Range("K20").Value = 3811.158
$ws.Range("M20").Value = -3564.158

# Row 94 (BSM)
$ws.Range("H94").Value = 2172.9412
$ws.Range("I94").Value = 1529.3334
$ws.Range("K94").Value = 1529.3334
$ws.Range("M94").Value = -1078.3334

# Row 134 (BSM)
$ws.Range("H134").Value = 5186.0713
$ws.Range("I134").Value = 4915.222
$ws.Range("J134").Value = 12499
$ws.Range("K134").Value = 14745.666
$ws.Range("L134").Value = 37497
$ws.Range("M134").Value = -12210.666
$ws.Range("N134").Value = -42567

$ws = $wb.Worksheets.Item("CRP")
# Row 16 (CRP)
$ws.Range("H16").Value = 3524.3157
$ws.Range("I16").Value = 1998.0714
$ws.Range("K16").Value = 1998.0714
$ws.Range("M16").Value = -1711.0714

# Row 22 (CRP)
$ws.Range("H22").Value = 553.6
$ws.Range("J22").Value = 1282.6666
$ws.Range("L22").Value = 1282.6666
$ws.Range("N22").Value = -1982.6666

# Row 60 (CRP)
$ws.Range("H60").Value = 19637.8
$ws.Range("I60").Value = 7549.5
$ws.Range("J60").Value = 27696.666
$ws.Range("K60").Value = 7549.5
$ws.Range("L60").Value = 27696.666
$ws.Range("M60").Value = -7038.5
$ws.Range("N60").Value = -28718.666

# Row 80 (CRP)
$ws.Range("H80").Value = 38888.46
$ws.Range("J80").Value = 38888.46
$ws.Range("L80").Value = 38888.46
$ws.Range("N80").Value = -41134.46

# Row 83 (CRP)
$ws.Range("H83").Value = 38888.46
$ws.Range("J83").Value = 38888.46
$ws.Range("L83").Value = 116665.38
$ws.Range("N83").Value = -127897.38

# Row 99 (CRP)
$ws.Range("H99").Value = 4333.6
$ws.Range("I99").Value = 4424.5
$ws.Range("K99").Value = 4424.5
$ws.Range("M99").Value = -2926.5

# Row 112 (CRP)
$ws.Range("H112").Value = 79973
$ws.Range("J112").Value = 79973
$ws.Range("L112").Value = 79973
$ws.Range("N112").Value = -82927

# Row 113 (CRP)
$ws.Range("H113").Value = 3524.3157
$ws.Range("I113").Value = 1998.0714
$ws.Range("K113").Value = 1998.0714
$ws.Range("M113").Value = 171.9286

# Row 126 (CRP)
$ws.Range("H126").Value = 4333.6
$ws.Range("I126").Value = 4424.5
$ws.Range("K126").Value = 13273.5
$ws.Range("M126").Value = -10803.5

# Row 132 (CRP)
$ws.Range("H132").Value = 3168.7083
$ws.Range("I132").Value = 3162.7273
$ws.Range("K132").Value = 9488.1819
$ws.Range("M132").Value = -6958.1819

# Row 134 (CRP)
$ws.Range("H134").Value = 7028.8823
$ws.Range("I134").Value = 5966.067
$ws.Range("J134").Value = 15000
$ws.Range("K134").Value = 17898.201
$ws.Range("L134").Value = 45000
$ws.Range("M134").Value = -15363.201
$ws.Range("N134").Value = -50070

$ws = $wb.Worksheets.Item("CUL")
# Row 34 (CUL)
$ws.Range("H34").Value = 357
$ws.Range("J34").Value = 750
$ws.Range("L34").Value = 2250
$ws.Range("N34").Value = -2418

# Row 68 (CUL)
$ws.Range("H68").Value = 3321.4285
$ws.Range("I68").Value = 3658.3333
$ws.Range("K68").Value = 10974.9999
$ws.Range("M68").Value = -10163.9999

# Row 71 (CUL)
$ws.Range("H71").Value = 3321.4285
$ws.Range("I71").Value = 3658.3333
$ws.Range("K71").Value = 32924.9997
$ws.Range("M71").Value = -28868.9997

# Row 113 (CUL)
$ws.Range("H113").Value = 1017.8
$ws.Range("J113").Value = 1072.5
$ws.Range("L113").Value = 3217.5
$ws.Range("N113").Value = -7557.5

$ws = $wb.Worksheets.Item("GSM")
# Row 32 (GSM)
$ws.Range("H32").Value = 36666.668
$ws.Range("J32").Value = 36666.668
$ws.Range("L32").Value = 36666.668
$ws.Range("N32").Value = -37258.668

# Row 97 (GSM)
$ws.Range("H97").Value = 7779.875
$ws.Range("I97").Value = 1647.8
$ws.Range("J97").Value = 18000
$ws.Range("K97").Value = 1647.8
$ws.Range("L97").Value = 18000
$ws.Range("M97").Value = -1151.8
$ws.Range("N97").Value = -18992

$ws = $wb.Worksheets.Item("LTW")
# Row 7 (LTW)
$ws.Range("H7").Value = 2514.5
$ws.Range("I7").Value = 2514.5
$ws.Range("K7").Value = 2514.5
$ws.Range("M7").Value = -2402.5

# Row 40 (LTW)
$ws.Range("H40").Value = 2421.7778
$ws.Range("J40").Value = 3099.5
$ws.Range("L40").Value = 3099.5
$ws.Range("N40").Value = -3371.5

# Row 59 (LTW)
$ws.Range("H59").Value = 200000
$ws.Range("J59").Value = 200000
$ws.Range("L59").Value = 200000

# Row 110 (LTW)
$ws.Range("H110").Value = 89900
$ws.Range("J110").Value = 89900
$ws.Range("L110").Value = 89900
$ws.Range("N110").Value = -98080

# Row 126 (LTW)
$ws.Range("H126").Value = 2514.5
$ws.Range("I126").Value = 2514.5
$ws.Range("K126").Value = 7543.5
$ws.Range("M126").Value = -5073.5

# Row 136 (LTW)
$ws.Range("H136").Value = 3543.4878
$ws.Range("I136").Value = 3189.2632
$ws.Range("J136").Value = 8030.3335
$ws.Range("K136").Value = 9567.7896
$ws.Range("L136").Value = 24091.0005
$ws.Range("M136").Value = -7017.7896
$ws.Range("N136").Value = -29191.0005

$ws = $wb.Worksheets.Item("WVR")
# Row 81 (WVR)
$ws.Range("H81").Value = 849.5
$ws.Range("I81").Value = 849.5
$ws.Range("K81").Value = 1699
$ws.Range("M81").Value = -638

# Row 84 (WVR)
$ws.Range("H84").Value = 849.5
$ws.Range("I84").Value = 849.5
$ws.Range("K84").Value = 8495
$ws.Range("M84").Value = -3191

# Row 103 (WVR)
$ws.Range("H103").Value = 33500
$ws.Range("J103").Value = 33500
$ws.Range("L103").Value = 33500
$ws.Range("N103").Value = -35844

# Row 126 (WVR)
$ws.Range("H126").Value = 3042.2307
$ws.Range("I126").Value = 2856.4211
$ws.Range("K126").Value = 8569.263300000001
$ws.Range("M126").Value = -6099.263300000001

# Row 132 (WVR)
$ws.Range("H132").Value = 4003.842
$ws.Range("I132").Value = 4168.0347
$ws.Range("J132").Value = 3474.7778
$ws.Range("K132").Value = 12504.1041
$ws.Range("L132").Value = 10424.3334
$ws.Range("M132").Value = -9974.1041
$ws.Range("N132").Value = -15484.3334

# Row 136 (WVR)
$ws.Range("H136").Value = 5377.522
$ws.Range("I136").Value = 3130.611
$ws.Range("K136").Value = 9391.832999999999
$ws.Range("M136").Value = -6841.832999999999

# Row 138 (WVR)
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()
